# Apply the "Add files via upload" edit to "Goals, QA, ASR and concerns.xlsx"
# Adds ASR (Architecturally Significant Requirement) rows/labels and a couple
# of new "Link with" cross-reference formulas on the Blad1 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- New ASR (Architecturally Significant Requirement) label rows ---
# (written first so the new shared-string table entries land in the same
# order as the source workbook: ASR-01..ASR-05 before the two descriptions)
$ws.Range("A56").Value = "ASR-01"
$ws.Range("A57").Value = "ASR-02"
$ws.Range("A58").Value = "ASR-03"
$ws.Range("A59").Value = "ASR-04"
$ws.Range("A60").Value = "ASR-05"

# --- New English translations for two existing concerns (C-01 / C-02) ---
$ws.Range("D36").Value = "Knowledge of systems and what is legally allowed "
$ws.Range("D37").Value = "Knowledge of encryption and standards is scarce"

# --- New / updated "Link with" formulas in column G (rows 48-53) ---
$ws.Range("G48").Formula = "=A31&"", ""&A44"
$ws.Range("G49").Formula = "=A28&"" ""&A27&"", ""&A36&"", ""&A37"
$ws.Range("G50").Formula = "=A27&"" ""&A45"
$ws.Range("G51").Formula = "=A31"
$ws.Range("G52").Formula = "=A29&"", ""&A38"
$ws.Range("G53").Formula = "=A27&"" ""&A39&"", ""&A40&"", ""&A41"

# --- Restore the selection to match the saved view state ---
$ws.Range("G53").Select()
